# Update the "CombineShipments" flag on the CreateNewCustomer sheet
# from "y" to "n" (cell H2).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateNewCustomer")
$ws.Range("H2").Value = "n"
